$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("vfiizu5d", "Training phase", 1, "['Purple', 'Orange', 'Orange', 'Orange', 'Green', 'Purple']", "[['', 'Red'], ['', 'Blue'], ['', 'Blue'], ['', 'Yellow'], ['', 'Blue'], ['', 'Yellow']]"),
    @("vfiizu5d", "Training phase", 2, "['Green', 'Green', 'Purple', 'Orange', 'Purple', 'Purple']", "[['', 'Red'], ['', 'Blue'], ['', 'Yellow'], ['', 'Red'], ['', 'Blue'], ['Yellow', '']]"),
    @("vfiizu5d", "Training phase", 3, "['Orange', 'Green', 'Purple', 'Orange', 'Orange', 'Green']", "[['', 'Red'], ['', 'Blue'], ['Yellow', ''], ['', 'Red'], ['', 'Blue'], ['', 'Blue']]"),
    @("vfiizu5d", "Training phase", 4, "['Orange', 'Purple', 'Orange', 'Purple', 'Green', 'Green']", "[['', 'Red'], ['', 'Blue'], ['', 'Yellow'], ['', 'Red'], ['', 'Blue'], ['', 'Red']]"),
    @("vfiizu5d", "Training phase", 5, "['Purple', 'Orange', 'Green', 'Green', 'Orange', 'Purple']", "[['', 'Red'], ['', 'Blue'], ['', 'Yellow'], ['', 'Blue'], ['', 'Red'], ['', 'Red']]"),
    @("vfiizu5d", "Test 1", 1, "['Green', 'Yellow', 'Purple', 'Red', 'Orange', 'Blue']", "[['', 'Red'], ['', 'Blue'], ['', 'Yellow'], ['', 'Red'], ['', 'Blue'], ['Yellow', '']]")
)

$startRow = 13
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
}
